$wb = $excel.ActiveWorkbook

# --- decision_var_bound: KLa LowerBound 0 -> 0.001 ---
$wsBound = $wb.Worksheets.Item("decision_var_bound")
$wsBound.Range("B6").Value = 0.001
$wsBound.Range("Q9").Select()

# --- decision_var: just move the selection cursor ---
$wsVar = $wb.Worksheets.Item("decision_var")
$wsVar.Range("A5").Select()

# --- fuzzy_goal: update goal targets, rename TSS goal to TN, add TKN goal row ---
$wsGoal = $wb.Worksheets.Item("fuzzy_goal")
$wsGoal.Range("B2").Value = 0
$wsGoal.Range("C2").Value = 100
$wsGoal.Range("B3").Value = 0
$wsGoal.Range("A5").Value = "Target_Effluent_TN (mg/L)"
$wsGoal.Range("B5").Value = 0
$wsGoal.Range("C5").Value = 10
$wsGoal.Range("B4").Value = 0
$wsGoal.Range("A6").Value = "Target_Effluent_TKN (mg/L)"
$wsGoal.Range("B6").Value = 0
$wsGoal.Range("C6").Value = 2
$wsGoal.Range("M17").Select()

# Restore the active sheet/tab selection to fuzzy_goal, as in the original workbook
$wsGoal.Activate()
